$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(0, 37.852828665146504, 32.085937424953364, 27.609281693794845, 43.465758266128674, 31.159768996379047, 31.196563288865786, 33.238953705806018),
    @(37.852828665146504, 0, 31.393557565474801, 82.031182071948862, 98.757015145876025, 67.054214637777434, 45.158836830976199, 75.30387532858056),
    @(32.085937424953364, 31.393557565474801, 0, 59.84370603723859, 89.361466028965737, 64.251603529056183, 50.661003479774394, 53.337692409187852),
    @(27.609281693794845, 82.031182071948862, 59.84370603723859, 0, 38.013632049814575, 48.068914723611122, 73.959778147718424, 41.944401626694649),
    @(43.465758266128674, 98.757015145876025, 89.361466028965737, 38.013632049814575, 0, 33.619985515373834, 63.130666499646544, 30.441279634974393),
    @(31.159768996379047, 67.054214637777434, 64.251603529056183, 48.068914723611122, 33.619985515373834, 0, 28.054480786641836, 35.393328591758291),
    @(31.196563288865786, 45.158836830976199, 50.661003479774394, 73.959778147718424, 63.130666499646544, 28.054480786641836, 0, 45.810508977011068),
    @(33.238953705806018, 75.30387532858056, 53.337692409187852, 41.944401626694649, 30.441279634974393, 35.393328591758291, 45.810508977011068, 0)
)

$rows = 8
$cols = 8
$data = New-Object 'object[,]' $rows,$cols
for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $data[$r,$c] = $values[$r][$c]
    }
}

$ws.Range("A1:H8").Value = $data
